$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Add the four new "Exterior / Frame" cap parts (#100-#103) that cover
# the exterior frame edges at the base and top of the printer.
# The cell-write order below reproduces the exact shared-string creation
# order captured in the target workbook (category first, then the two
# base-cap filenames, then the base-cap names written out of row order,
# then the top-cap rows filled name-then-filename).
# -------------------------------------------------------------------

# Row 71 - #100 Base Cap Type 1
$ws.Range("A71").Value = 100
$ws.Range("B71").Value = "Exterior"
$ws.Range("C71").Value = "Frame"
$ws.Range("D71").Value = "N"
$ws.Range("F71").Value = "ABS"
$ws.Range("G71").Value = 2

# Row 72 - #101 Base Cap Type 2
$ws.Range("A72").Value = 101
$ws.Range("B72").Value = "Exterior"
$ws.Range("C72").Value = "Frame"
$ws.Range("D72").Value = "N"
$ws.Range("F72").Value = "ABS"
$ws.Range("G72").Value = 2

$ws.Range("I71").Value = "100 - Exterior - Frame - Base Cap Type 1.stl"
$ws.Range("I72").Value = "101 - Exterior - Frame - Base Cap Type 2.stl"
$ws.Range("E72").Value = "Base Cap Type 2"
$ws.Range("E71").Value = "Base Cap Type 1"

# Row 73 - #102 Top Cap Type 1
$ws.Range("A73").Value = 102
$ws.Range("B73").Value = "Exterior"
$ws.Range("C73").Value = "Frame"
$ws.Range("D73").Value = "N"
$ws.Range("F73").Value = "ABS"
$ws.Range("G73").Value = 2
$ws.Range("E73").Value = "Top Cap Type 1"
$ws.Range("I73").Value = "102 - Exterior - Frame - Top Cap Type 1.stl"

# Row 74 - #103 Top Cap Type 2
$ws.Range("A74").Value = 103
$ws.Range("B74").Value = "Exterior"
$ws.Range("C74").Value = "Frame"
$ws.Range("D74").Value = "N"
$ws.Range("F74").Value = "ABS"
$ws.Range("G74").Value = 2
$ws.Range("E74").Value = "Top Cap Type 2"
$ws.Range("I74").Value = "103 - Exterior - Frame - Top Cap Type 2.stl"

# Grow Table1 to cover the four new rows.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:I74"))

# Keep the sheet view pinned near the newly-added rows, matching the
# author finishing their edit at the bottom of the table.
[void]$ws.Range("I74").Select()
$excel.ActiveWindow.ScrollRow = 47
